# Insert a new price-report row for "Macroferia Regional de Talca" / Frutilla
# at sheet row 740, pushing the existing rows 740-796 down to 741-797.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 740 (shifts 740:796 -> 741:797, keeps formats).
$ws.Rows.Item(740).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A740").Value = 5
$ws.Range("B740").Value = "Macroferia Regional de Talca"
$ws.Range("C740").Value = "Maule"
$ws.Range("D740").Value = 45021
$ws.Range("E740").Value = 7
$ws.Range("F740").Value = "Fruta"
$ws.Range("G740").Value = 100101
$ws.Range("H740").Value = "Berries"
$ws.Range("I740").Value = 100112025
$ws.Range("J740").Value = "Frutilla"
$ws.Range("K740").Value = "Sin especificar"
$ws.Range("L740").Value = "Especial"
$ws.Range("M740").Value = 200
$ws.Range("N740").Value = 9000
$ws.Range("O740").Value = 9000
$ws.Range("P740").Value = 9000
$ws.Range("Q740").Value = "$/caja 7 kilos"
$ws.Range("R740").Value = "Región del Maule"
$ws.Range("S740").Value = 1286
$ws.Range("T740").Value = 7
